$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" header suffixes to "_FV2404" / "_FV2410"
$fields = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fields[$i] + "_FV2404"
}
for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fields[$i] + "_FV2410"
}

# Turn the data range into a proper Excel table ("Table1")
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
